# Auto-generated script applying the cryptos.xlsx price/volume/coin updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue {
    param($Cell, [string]$Text)
    # Prefix with an apostrophe so Excel stores the exact literal text
    # (keeps "1.00"-style values as strings instead of auto-converting to numbers)
    $Cell.Value = "'" + $Text
}

Set-TextValue $ws.Range("D2") "65.467.25"
$ws.Range("E2").Value = "  +1.50%  "

Set-TextValue $ws.Range("D3") "3.155.90"
$ws.Range("E3").Value = "  +4.10%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.11%  "

Set-TextValue $ws.Range("D5") "568.56"
$ws.Range("E5").Value = "  +1.48%  "

Set-TextValue $ws.Range("D6") "151.68"
$ws.Range("E6").Value = "  +8.44%  "

Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.05%  "

Set-TextValue $ws.Range("D8") "3.149.31"
$ws.Range("E8").Value = "  +4.27%  "

Set-TextValue $ws.Range("D9") "0.505"
$ws.Range("E9").Value = "  +5.27%  "

Set-TextValue $ws.Range("D10") "6.79"
$ws.Range("E10").Value = "  +12.80%  "

Set-TextValue $ws.Range("D11") "0.163"
$ws.Range("E11").Value = "  +5.76%  "

Set-TextValue $ws.Range("D12") "0.470"
$ws.Range("E12").Value = "  +4.75%  "

Set-TextValue $ws.Range("D13") "37.68"
$ws.Range("E13").Value = "  +10.56%  "

Set-TextValue $ws.Range("D14") "0.0000227"
$ws.Range("E14").Value = "  +6.87%  "

Set-TextValue $ws.Range("D15") "3.666.82"
$ws.Range("E15").Value = "  +4.17%  "

Set-TextValue $ws.Range("D16") "65.352.61"
$ws.Range("E16").Value = "  +1.43%  "

Set-TextValue $ws.Range("D17") "3.159.90"
$ws.Range("E17").Value = "  +4.21%  "

Set-TextValue $ws.Range("D18") "0.113"
$ws.Range("E18").Value = "  +2.11%  "

Set-TextValue $ws.Range("D19") "520.78"
$ws.Range("E19").Value = "  +7.96%  "

Set-TextValue $ws.Range("D20") "6.89"
$ws.Range("E20").Value = "  +6.81%  "

Set-TextValue $ws.Range("D21") "14.17"
$ws.Range("E21").Value = "  +6.17%  "

Set-TextValue $ws.Range("D22") "0.714"
$ws.Range("E22").Value = "  +8.79%  "

Set-TextValue $ws.Range("D23") "7.52"
$ws.Range("E23").Value = "  +8.53%  "

Set-TextValue $ws.Range("D24") "12.98"
$ws.Range("E24").Value = "  +4.40%  "

Set-TextValue $ws.Range("D25") "79.93"
$ws.Range("E25").Value = "  +2.85%  "

Set-TextValue $ws.Range("D26") "0.998"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D27") "8.95"
$ws.Range("E27").Value = "  +19.45%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D28") "2.84"
$ws.Range("E28").Value = "  +4.71%  "

Set-TextValue $ws.Range("D29") "2.18"
$ws.Range("E29").Value = "  +5.93%  "

$ws.Range("B30").Value = "Stacks"
$ws.Range("C30").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D30") "2.76"
$ws.Range("E30").Value = "  +7.96%  "

$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D31") "1.00"
$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D32") "27.01"
$ws.Range("E32").Value = "  +5.66%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D33") "588.29"
$ws.Range("E33").Value = "  +13.08%  "

Set-TextValue $ws.Range("D34") "1.16"
$ws.Range("E34").Value = "  +4.41%  "

Set-TextValue $ws.Range("D35") "5.68"
$ws.Range("E35").Value = "  +8.37%  "

Set-TextValue $ws.Range("D36") "6.20"
$ws.Range("E36").Value = "  +7.85%  "

Set-TextValue $ws.Range("D37") "53.43"
$ws.Range("E37").Value = "  +3.10%  "

Set-TextValue $ws.Range("D38") "0.0429"
$ws.Range("E38").Value = "  +5.82%  "

Set-TextValue $ws.Range("D39") "3.06"
$ws.Range("E39").Value = "  +27.25%  "

Set-TextValue $ws.Range("D40") "0.0834"
$ws.Range("E40").Value = "  +5.76%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D41") "0.123"
$ws.Range("E41").Value = "  +3.53%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D42") "3.072.47"
$ws.Range("E42").Value = "  +9.44%  "

Set-TextValue $ws.Range("D43") "8.40"
$ws.Range("E43").Value = "  +3.85%  "

Set-TextValue $ws.Range("D44") "0.260"
$ws.Range("E44").Value = "  +9.87%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D45") "2.23"
$ws.Range("E45").Value = "  +13.15%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D46") "0.999"
$ws.Range("E46").Value = "  +0.04%  "

Set-TextValue $ws.Range("D47") "26.18"
$ws.Range("E47").Value = "  +11.81%  "

Set-TextValue $ws.Range("D48") "0.0₃0551"
$ws.Range("E48").Value = "  +8.30%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D49") "121.95"
$ws.Range("E49").Value = "  +5.17%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D50") "0.111"
$ws.Range("E50").Value = "  +4.80%  "

Set-TextValue $ws.Range("D51") "2.17"
$ws.Range("E51").Value = "  +9.11%  "

